$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column widths: column B becomes narrower, column C gets an explicit width
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Columns.Item(3).ColumnWidth = 17.375

# Fill column N with formulas building an INSERT SQL statement from columns A, B, C.
# N2 is written alone (non-shared formula). N3:N37 are written together as a
# single range-formula assignment so Excel stores it as a shared formula group
# (matches the xr:si shared-formula layout Excel itself produces on fill).
$ws.Range("N2").Formula = '="INSERT INTO `stores` (`name`, `address`, `image`) VALUES (''"&A2&"'',''"&B2&"'',''"&C2&"'');"'
$ws.Range("N3:N37").Formula = '="INSERT INTO `stores` (`name`, `address`, `image`) VALUES (''"&A3&"'',''"&B3&"'',''"&C3&"'');"'

# Update view: select N2:N37 with active cell N2 (also clears the stale
# topLeftCell scroll position left over from the previous selection)
$ws.Range("N2:N37").Select() | Out-Null
